# final ve üçüncülük maçı görüntüleri eklendi, 3.lük maçı skoru düzeltildi
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final maçı (row 21) için video linki ekle
$ws.Range("K21").Value = "https://youtu.be/UwPkIRGWSeA"

# 3.lük maçı (row 24) skoru düzeltildi: D24 4 -> 3
$ws.Range("D24").Value = 3

# 3.lük maçı (row 24) için video linki ekle
$ws.Range("K24").Value = "https://youtu.be/v4QNOTOEC-E"

# Seçili hücreyi güncelle
$ws.Range("K27").Select()
